$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2214.4285
$ws.Range("I43").Value = 750.5
$ws.Range("J43").Value = 2800
$ws.Range("K43").Value = 750.5
$ws.Range("L43").Value = 2800
$ws.Range("M43").Value = -681.5
$ws.Range("N43").Value = -2938
$ws.Range("H55").Value = 390
$ws.Range("I55").Value = 43.333332
$ws.Range("J55").Value = 650
$ws.Range("K55").Value = 43.333332
$ws.Range("L55").Value = 650
$ws.Range("M55").Value = 170.666668
$ws.Range("N55").Value = -1078
$ws.Range("H88").Value = 618961.9399999999
$ws.Range("I88").Value = 457.9091
$ws.Range("J88").Value = 1374911.4
$ws.Range("K88").Value = 457.9091
$ws.Range("L88").Value = 1374911.4
$ws.Range("M88").Value = -51.90910000000002
$ws.Range("N88").Value = -1375723.4
$ws.Range("H91").Value = 618961.9399999999
$ws.Range("I91").Value = 457.9091
$ws.Range("J91").Value = 1374911.4
$ws.Range("K91").Value = 457.9091
$ws.Range("L91").Value = 1374911.4
$ws.Range("M91").Value = 946.0908999999999
$ws.Range("N91").Value = -1377719.4
$ws.Range("H137").Value = 1377.4688
$ws.Range("I137").Value = 1178.0233
$ws.Range("J137").Value = 1785.8572
$ws.Range("K137").Value = 3534.0699
$ws.Range("L137").Value = 5357.571599999999
$ws.Range("M137").Value = -984.0699000000004
$ws.Range("N137").Value = -10457.5716

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 803.5
$ws.Range("I4").Value = 500
$ws.Range("K4").Value = 500
$ws.Range("M4").Value = -384
$ws.Range("H5").Value = 255
$ws.Range("I5").Value = 184
$ws.Range("J5").Value = 373.33334
$ws.Range("K5").Value = 184
$ws.Range("L5").Value = 373.33334
$ws.Range("M5").Value = -72
$ws.Range("N5").Value = -597.33334
$ws.Range("H32").Value = 6511.11
$ws.Range("I32").Value = 4821.1777
$ws.Range("J32").Value = 21720.5
$ws.Range("K32").Value = 4821.1777
$ws.Range("L32").Value = 21720.5
$ws.Range("M32").Value = -4534.1777
$ws.Range("N32").Value = -22294.5
$ws.Range("H43").Value = 8688.5
$ws.Range("J43").Value = 8688.5
$ws.Range("L43").Value = 8688.5
$ws.Range("N43").Value = -9314.5
$ws.Range("H61").Value = 90910840
$ws.Range("I61").Value = 142858340
$ws.Range("K61").Value = 142858340
$ws.Range("M61").Value = -142858128
$ws.Range("H74").Value = 3001.6191
$ws.Range("I74").Value = 2440.3
$ws.Range("K74").Value = 2440.3
$ws.Range("M74").Value = -1566.3
$ws.Range("H77").Value = 3001.6191
$ws.Range("I77").Value = 2440.3
$ws.Range("K77").Value = 12201.5
$ws.Range("M77").Value = -7833.5
$ws.Range("H132").Value = 3781.6785
$ws.Range("I132").Value = 3567.4167
$ws.Range("K132").Value = 10702.2501
$ws.Range("M132").Value = -8172.250100000001
$ws.Range("H136").Value = 90910840
$ws.Range("I136").Value = 142858340
$ws.Range("K136").Value = 428575020
$ws.Range("M136").Value = -428572470

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 255
$ws.Range("I4").Value = 184
$ws.Range("J4").Value = 373.33334
$ws.Range("K4").Value = 184
$ws.Range("L4").Value = 373.33334
$ws.Range("M4").Value = -69
$ws.Range("N4").Value = -603.33334
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1196.9385
$ws.Range("I31").Value = 1152
$ws.Range("J31").Value = 1569.2858
$ws.Range("K31").Value = 1152
$ws.Range("L31").Value = 1569.2858
$ws.Range("M31").Value = -857
$ws.Range("N31").Value = -2159.2858
$ws.Range("H34").Value = 1196.9385
$ws.Range("I34").Value = 1152
$ws.Range("J34").Value = 1569.2858
$ws.Range("K34").Value = 1152
$ws.Range("L34").Value = 1569.2858
$ws.Range("M34").Value = -950
$ws.Range("N34").Value = -1973.2858
$ws.Range("H52").Value = 27801.285
$ws.Range("I52").Value = 14854.5
$ws.Range("J52").Value = 32980
$ws.Range("K52").Value = 14854.5
$ws.Range("L52").Value = 32980
$ws.Range("M52").Value = -14560.5
$ws.Range("N52").Value = -33568
$ws.Range("H86").Value = 5641640
$ws.Range("I86").Value = 7488518
$ws.Range("J86").Value = 101005.664
$ws.Range("K86").Value = 7488518
$ws.Range("L86").Value = 101005.664
$ws.Range("M86").Value = -7487395
$ws.Range("N86").Value = -103251.664
$ws.Range("H89").Value = 5641640
$ws.Range("I89").Value = 7488518
$ws.Range("J89").Value = 101005.664
$ws.Range("K89").Value = 37442590
$ws.Range("L89").Value = 505028.32
$ws.Range("M89").Value = -37436974
$ws.Range("N89").Value = -516260.32

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5750
$ws.Range("J80").Value = 5750
$ws.Range("L80").Value = 17250
$ws.Range("N80").Value = -19122
$ws.Range("H83").Value = 5750
$ws.Range("J83").Value = 5750
$ws.Range("L83").Value = 51750
$ws.Range("N83").Value = -61110
$ws.Range("H138").Value = 2654.8684
$ws.Range("I138").Value = 2692.2307
$ws.Range("J138").Value = 2635.44
$ws.Range("K138").Value = 8076.6921
$ws.Range("L138").Value = 7906.32
$ws.Range("M138").Value = -2936.6921
$ws.Range("N138").Value = -18186.32
$ws.Range("H140").Value = 23653.188
$ws.Range("J140").Value = 3480.4814
$ws.Range("L140").Value = 10441.4442
$ws.Range("N140").Value = -20801.4442

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20576
$ws.Range("H80").Value = 6350
$ws.Range("I80").Value = 5000
$ws.Range("K80").Value = 5000
$ws.Range("M80").Value = -4002
$ws.Range("H81").Value = 20000
$ws.Range("J81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -21996
$ws.Range("H83").Value = 6350
$ws.Range("I83").Value = 5000
$ws.Range("K83").Value = 25000
$ws.Range("M83").Value = -20008
$ws.Range("H84").Value = 20000
$ws.Range("J84").Value = 20000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -69984
$ws.Range("H136").Value = 12074.728
$ws.Range("J136").Value = 12074.728
$ws.Range("L136").Value = 36224.18399999999
$ws.Range("N136").Value = -41324.18399999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2853.9
$ws.Range("I7").Value = 2591.75
$ws.Range("J7").Value = 3902.5
$ws.Range("K7").Value = 2591.75
$ws.Range("L7").Value = 3902.5
$ws.Range("M7").Value = -2479.75
$ws.Range("N7").Value = -4126.5
$ws.Range("H22").Value = 1402.75
$ws.Range("I22").Value = 1055.5
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 1055.5
$ws.Range("L22").Value = 1750
$ws.Range("M22").Value = -760.5
$ws.Range("N22").Value = -2340
$ws.Range("H27").Value = 1402.75
$ws.Range("I27").Value = 1055.5
$ws.Range("J27").Value = 1750
$ws.Range("K27").Value = 1055.5
$ws.Range("L27").Value = 1750
$ws.Range("M27").Value = -948.5
$ws.Range("N27").Value = -1964
$ws.Range("H61").Value = 1237.7858
$ws.Range("J61").Value = 1417.2
$ws.Range("L61").Value = 1417.2
$ws.Range("N61").Value = -1821.2
$ws.Range("H93").Value = 1420.6
$ws.Range("I93").Value = 1400.6666
$ws.Range("K93").Value = 1400.6666
$ws.Range("M93").Value = -152.6666
$ws.Range("H109").Value = 25285
$ws.Range("J109").Value = 25285
$ws.Range("L109").Value = 25285
$ws.Range("N109").Value = -28059
$ws.Range("H113").Value = 1237.7858
$ws.Range("J113").Value = 1417.2
$ws.Range("L113").Value = 1417.2
$ws.Range("N113").Value = -5757.2
$ws.Range("H126").Value = 2853.9
$ws.Range("I126").Value = 2591.75
$ws.Range("J126").Value = 3902.5
$ws.Range("K126").Value = 7775.25
$ws.Range("L126").Value = 11707.5
$ws.Range("M126").Value = -5305.25
$ws.Range("N126").Value = -16647.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 65003.5
$ws.Range("J12").Value = 65003.5
$ws.Range("L12").Value = 65003.5
$ws.Range("N12").Value = -65287.5
$ws.Range("H96").Value = 10000
$ws.Range("I96").Value = 10000
$ws.Range("K96").Value = 10000
$ws.Range("M96").Value = -8627
$ws.Range("H100").Value = 460
$ws.Range("I100").Value = 460
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 920
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -379
$ws.Range("N100").Value = ""
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800
$ws.Range("H132").Value = 2247.0293
$ws.Range("I132").Value = 2034.1111
$ws.Range("K132").Value = 6102.3333
$ws.Range("M132").Value = -3572.3333
